# Apply the crypto price/volume refresh captured by the commit
# "Updated cryptos list on Sat May 20 03:12:49 UTC 2023 with GitHub Actions".
#
# Column D ("Price") holds numeric-looking text (e.g. "26.902.25", "309.39")
# that must stay plain text, exactly as scraped - it is not a real number
# (some prices even use two dots as thousands+decimal separators). Excel
# auto-converts bare numeric-looking strings assigned via .Value into real
# numbers, so a leading apostrophe is used to force text entry, matching the
# original inline-string cell type.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '''26.902.25'
$ws.Range('E2').Value = '  +0.33%  '

# Row 3
$ws.Range('D3').Value = '''1.813.16'
$ws.Range('E3').Value = '  +0.73%  '

# Row 4
$ws.Range('E4').Value = '  +0.13%  '

# Row 5
$ws.Range('D5').Value = '''309.39'

# Row 6
$ws.Range('E6').Value = '  +0.11%  '

# Row 7
$ws.Range('D7').Value = '''0.4666'
$ws.Range('E7').Value = '  +0.55%  '

# Row 8
$ws.Range('D8').Value = '''0.3656'
$ws.Range('E8').Value = '  -0.76%  '

# Row 9
$ws.Range('D9').Value = '''0.07351'
$ws.Range('E9').Value = '  +0.02%  '

# Row 10
$ws.Range('D10').Value = '''0.8686'
$ws.Range('E10').Value = '  +0.31%  '

# Row 11
$ws.Range('D11').Value = '''20.31'
$ws.Range('E11').Value = '  -0.21%  '

# Row 12
$ws.Range('D12').Value = '''1.844.78'
$ws.Range('E12').Value = '  -2.66%  '

# Row 13
$ws.Range('D13').Value = '''5.385'
$ws.Range('E13').Value = '  +0.74%  '

# Row 14
$ws.Range('D14').Value = '''0.07088'
$ws.Range('E14').Value = '  +0.87%  '

# Row 15
$ws.Range('D15').Value = '''6.512'
$ws.Range('E15').Value = '  -0.03%  '

# Row 16
$ws.Range('D16').Value = '''91.62'
$ws.Range('E16').Value = '  +0.32%  '

# Row 17
$ws.Range('E17').Value = '  +0.13%  '

# Row 18
$ws.Range('D18').Value = '''0.000008705'
$ws.Range('E18').Value = '  +0.35%  '

# Row 19
$ws.Range('E19').Value = '  +0.07%  '

# Row 20
$ws.Range('D20').Value = '''14.64'
$ws.Range('E20').Value = '  +0.19%  '

# Row 21
$ws.Range('D21').Value = '''26.915.38'
$ws.Range('E21').Value = '  +0.36%  '

# Row 22
$ws.Range('D22').Value = '''5.299'
$ws.Range('E22').Value = '  +0.21%  '

# Row 23
$ws.Range('D23').Value = '''10.63'
$ws.Range('E23').Value = '  +0.58%  '

# Row 24
$ws.Range('D24').Value = '''2.033.93'
$ws.Range('E24').Value = '  -3.71%  '

# Row 25
$ws.Range('E25').Value = '  -0.53%  '

# Row 26
$ws.Range('D26').Value = '''150.75'
$ws.Range('E26').Value = '  -0.38%  '

# Row 27
$ws.Range('D27').Value = '''18.33'
$ws.Range('E27').Value = '  +0.16%  '

# Row 28
$ws.Range('D28').Value = '''2.147'
$ws.Range('E28').Value = '  +1.17%  '

# Row 29
$ws.Range('D29').Value = '''5.266'
$ws.Range('E29').Value = '  +0.73%  '

# Row 30
$ws.Range('D30').Value = '''115.26'
$ws.Range('E30').Value = '  -0.36%  '

# Row 31
$ws.Range('D31').Value = '''0.08931'
$ws.Range('E31').Value = '  +0.42%  '

# Row 32
$ws.Range('D32').Value = '''0.7538'
$ws.Range('E32').Value = '  -0.13%  '

# Row 33
$ws.Range('D33').Value = '''1.155'
$ws.Range('E33').Value = '  +0.95%  '

# Row 34
$ws.Range('D34').Value = '''4.489'
$ws.Range('E34').Value = '  +1.08%  '

# Row 35
$ws.Range('E35').Value = '  -0.73%  '

# Row 37
$ws.Range('D37').Value = '''1.084'
$ws.Range('E37').Value = '  -1.70%  '

# Row 38
$ws.Range('D38').Value = '''0.05281'
$ws.Range('E38').Value = '  +0.83%  '

# Row 39
$ws.Range('D39').Value = '''0.01948'
$ws.Range('E39').Value = '  -0.01%  '

# Row 40
$ws.Range('D40').Value = '''2.973'
$ws.Range('E40').Value = '  +1.42%  '

# Row 41
$ws.Range('D41').Value = '''7.210'
$ws.Range('E41').Value = '  +0.32%  '

# Row 42
$ws.Range('D42').Value = '''0.5300'
$ws.Range('E42').Value = '  +0.66%  '

# Row 43
$ws.Range('D43').Value = '''2.280'
$ws.Range('E43').Value = '  -2.81%  '

# Row 44
$ws.Range('D44').Value = '''0.1653'
$ws.Range('E44').Value = '  -0.16%  '

# Row 45
$ws.Range('D45').Value = '''8.401'
$ws.Range('E45').Value = '  -0.70%  '

# Row 46
$ws.Range('D46').Value = '''0.4872'
$ws.Range('E46').Value = '  -2.40%  '

# Row 47
$ws.Range('D47').Value = '''10.44'
$ws.Range('E47').Value = '  +1.79%  '

# Row 48
$ws.Range('E48').Value = '  +0.13%  '

# Rows 49-50: NEARProtocol and Quant swapped ranking positions in this
# refresh, so row 49 now shows Quant (with its updated price/volume) and
# row 50 shows NEARProtocol (with its updated volume).
# Row 49
$ws.Range('B49').Value = 'Quant'
$ws.Range('C49').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D49').Value = '''103.06'
$ws.Range('E49').Value = '  -0.90%  '

# Row 50
$ws.Range('B50').Value = 'NEARProtocol'
$ws.Range('C50').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D50').Value = '''1.657'
$ws.Range('E50').Value = '  -0.17%  '

# Row 51
$ws.Range('E51').Value = '  +0.05%  '
